$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation arrived for "Terminal La Palmera de La Serena -
# Espinaca": insert it as the new row 143, pushing the existing rows
# 143-146 down to 144-147 (the data itself is unchanged, only its row
# position shifts).
$ws.Rows(143).Insert()

$ws.Cells.Item(143, 1).Value = 8
$ws.Cells.Item(143, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(143, 3).Value = "Coquimbo"
$ws.Cells.Item(143, 4).Value = 44448
$ws.Cells.Item(143, 5).Value = 4
$ws.Cells.Item(143, 6).Value = 100112012
$ws.Cells.Item(143, 7).Value = "Espinaca"
$ws.Cells.Item(143, 8).Value = "Sin especificar"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 3000
$ws.Cells.Item(143, 11).Value = 400
$ws.Cells.Item(143, 12).Value = 500
$ws.Cells.Item(143, 13).Value = 450
$ws.Cells.Item(143, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(143, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(143, 16).Value = 900
$ws.Cells.Item(143, 17).Value = 0.5
$ws.Cells.Item(143, 18).Value = "Hortaliza"
